# Add a new "Participation and Activities" column (G) to the grade sheet,
# populate it with a flat participation score of 5 for every student row
# (rows 2-35), size the new column, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header in G1
$ws.Range("G1").Value = "Participation and Activities"

# Participation score of 5 for each of the 34 student rows (2 through 35)
for ($r = 2; $r -le 35; $r++) {
    $ws.Cells.Item($r, 7).Value = 5
}

# Widen column G to fit the new header text
$ws.Columns.Item(7).ColumnWidth = 20

# Move the selection to I19 (clears the old A33 top-left / F37:F38 selection)
[void]$ws.Range("I19").Select()
